$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Row 207: add the "Grant / 2017 / Moonshot / DCTD" tag columns (B:E) ---
# That 4-value pattern already exists all over the sheet (e.g. B2:E2); reuse
# it here and paint it with this row's own highlight colour (yellow).
$ws.Range("B2:E2").Copy() | Out-Null
$ws.Range("B207:E207").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$ws.Range("B207:E207").NumberFormat = "@"
$ws.Range("B207:E207").Interior.Color = $ws.Cells.Item(207, 1).Interior.Color

# --- 2) Row 427: same tag columns, painted with this row's highlight colour (red) ---
$ws.Range("B2:E2").Copy() | Out-Null
$ws.Range("B427:E427").PasteSpecial(-4163) | Out-Null
$ws.Range("B427:E427").NumberFormat = "@"
$ws.Range("B427:E427").Interior.Color = $ws.Cells.Item(427, 1).Interior.Color

# --- 3) Insert a brand-new row for a missed "single result" DOM-parse test
#        case, with two plain (unstyled) cells ---
$ws.Rows.Item(428).Insert()
$ws.Rows.Item(428).ClearFormats() | Out-Null
$ws.Range("B428:E428").ClearContents() | Out-Null
$ws.Cells.Item(428, 1).Value2 = "asdfasdfasdf"
$ws.Cells.Item(428, 6).Value2 = "No results"

# --- 4) Row 430 (was 429 before the insert): tag columns with no fill ---
$ws.Range("B2:E2").Copy() | Out-Null
$ws.Range("B430:E430").PasteSpecial(-4163) | Out-Null
$ws.Range("B430:E430").NumberFormat = "@"

# --- 5) Row 438 (was 437 before the insert): flag it with the same
#        "Moonshot of interest that break" note used on row 427 ---
$ws.Range("F427").Copy() | Out-Null
$ws.Range("F438").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(438, 6).Value2 = $ws.Cells.Item(427, 6).Value2

# --- 6) Selection / viewport bookkeeping ---
$ws.Range("C438").Select() | Out-Null
